$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.445.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.08%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.939.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.12%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  -0.82%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.49%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'56.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.12%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -4.36%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0808"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.40%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.09%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.222.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.32%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'21.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.23%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.802"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.61%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'13.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.48%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'5.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.49%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.931.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.18%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'36.392.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.02%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'68.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.10%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0852"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.59%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'226.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.05%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.28%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.18%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -5.65%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.11%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -5.01%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'159.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "'  +7.56%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.42%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.16%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D32").Value = "'4.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.46%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -4.73%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.54%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.13%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -1.34%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -1.29%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.49%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +9.82%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.09%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.02%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.71%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -4.74%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'15.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.19%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -3.07%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.328.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.53%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'85.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.05%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.89%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.63%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +14.85%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.115.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.31%  "
$ws.Range("E51").Style = "Normal"
